$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column D (截止一直未充电时间) for all data rows 2-58 to the new refresh timestamp
$ws.Range("D2:D58").Value = 45953.287164351852

# Rows 18-53: refreshed charging-terminal data (new A/B/C values)
$ws.Range("A18").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B18").Value = "603号直流"
$ws.Range("C18").Value = 45950.523645833331
$ws.Range("A19").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B19").Value = "101号直流"
$ws.Range("C19").Value = 45950.574606481481
$ws.Range("A20").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B20").Value = "501号直流"
$ws.Range("C20").Value = 45950.736122685186
$ws.Range("A21").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B21").Value = "602号直流"
$ws.Range("C21").Value = 45951.042326388888
$ws.Range("A22").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B22").Value = "502号直流"
$ws.Range("C22").Value = 45951.226111111115
$ws.Range("A23").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B23").Value = "103号直流"
$ws.Range("C23").Value = 45951.540775462963
$ws.Range("A24").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B24").Value = "102号直流"
$ws.Range("C24").Value = 45951.62295138889
$ws.Range("A25").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B25").Value = "104号直流"
$ws.Range("C25").Value = 45952.059756944444
$ws.Range("A26").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B26").Value = "306号直流"
$ws.Range("C26").Value = 45952.066192129627
$ws.Range("A27").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B27").Value = "401号直流"
$ws.Range("C27").Value = 45952.110937500001
$ws.Range("A28").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B28").Value = "904号直流"
$ws.Range("C28").Value = 45952.248090277775
$ws.Range("A29").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B29").Value = "702号直流"
$ws.Range("C29").Value = 45952.248715277776
$ws.Range("A30").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B30").Value = "901号直流"
$ws.Range("C30").Value = 45952.250949074078
$ws.Range("A31").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B31").Value = "210号直流"
$ws.Range("C31").Value = 45952.303078703706
$ws.Range("A32").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B32").Value = "406号直流"
$ws.Range("C32").Value = 45952.398726851854
$ws.Range("A33").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B33").Value = "A01号直流"
$ws.Range("C33").Value = 45952.478333333333
$ws.Range("A34").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B34").Value = "201号直流"
$ws.Range("C34").Value = 45952.525138888886
$ws.Range("A35").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B35").Value = "102号直流"
$ws.Range("C35").Value = 45952.540023148147
$ws.Range("A36").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B36").Value = "104号直流"
$ws.Range("C36").Value = 45952.551782407405
$ws.Range("A37").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B37").Value = "903号直流"
$ws.Range("C37").Value = 45952.557546296295
$ws.Range("A38").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B38").Value = "110号直流"
$ws.Range("C38").Value = 45952.564340277779
$ws.Range("A39").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B39").Value = "106号直流"
$ws.Range("C39").Value = 45952.571180555555
$ws.Range("A40").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B40").Value = "503号直流"
$ws.Range("C40").Value = 45952.577430555553
$ws.Range("A41").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B41").Value = "403号直流"
$ws.Range("C41").Value = 45952.582638888889
$ws.Range("A42").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B42").Value = "501号直流"
$ws.Range("C42").Value = 45952.585428240738
$ws.Range("A43").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B43").Value = "312号直流"
$ws.Range("C43").Value = 45952.625532407408
$ws.Range("A44").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B44").Value = "301号直流"
$ws.Range("C44").Value = 45952.627812500003
$ws.Range("A45").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B45").Value = "304号直流"
$ws.Range("C45").Value = 45952.653726851851
$ws.Range("A46").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B46").Value = "904号直流"
$ws.Range("C46").Value = 45952.658252314817
$ws.Range("A47").Value = "长沙特来电飞狐四方坪西区充电站"
$ws.Range("B47").Value = "B04号直流"
$ws.Range("C47").Value = 45952.670949074076
$ws.Range("A48").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B48").Value = "202号直流"
$ws.Range("C48").Value = 45952.672708333332
$ws.Range("A49").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B49").Value = "208号直流"
$ws.Range("C49").Value = 45952.681574074071
$ws.Range("A50").Value = "长沙特来电飞狐四方坪东区充电站"
$ws.Range("B50").Value = "101号直流"
$ws.Range("C50").Value = 45952.68204861111
$ws.Range("A51").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B51").Value = "206号直流"
$ws.Range("C51").Value = 45952.69321759259
$ws.Range("A52").Value = "长沙市开福区高岭香江国际城充电站建设项目"
$ws.Range("B52").Value = "111号直流"
$ws.Range("C52").Value = 45952.695925925924
$ws.Range("A53").Value = "长沙特来电飞狐四方坪南区充电站"
$ws.Range("B53").Value = "202号直流"
$ws.Range("C53").Value = 45952.780150462961

# Rows 54-58: no longer have data -- clear content but keep formatting
$ws.Range("A54:D54").ClearContents()
$ws.Range("A55:D55").ClearContents()
$ws.Range("A56:D56").ClearContents()
$ws.Range("A57:D57").ClearContents()
$ws.Range("A58:D58").ClearContents()

# Restore the active-cell selection recorded in the sheet view
$ws.Range("I18").Select()
